# Add a new "2022-Q1" sheet (same layout as the other quarterly sheets)
# right before the "总计" (totals) sheet, and add a corresponding summary
# row at the top of the "总计" sheet's data.

$wb = $excel.ActiveWorkbook

# Locate the "总计" sheet (last sheet in the workbook) and the most recent
# quarterly sheet (2021-Q4), whose layout the new sheet will copy.
$srcSheet = $wb.Worksheets.Item(3)
$totalSheet = $wb.Worksheets.Item("总计")

# Duplicate 2021-Q4 so the new sheet inherits the same sheetPr/column
# layout/header style, placing the copy immediately before "总计". (Note:
# re-fetch "总计" by name afterwards rather than reusing $totalSheet.Index
# — the cached Index on the old reference does not update after the
# insert shifts it over.)
$srcSheet.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item("总计").Previous
$newSheet.Name = "2022-Q1"

# Drop the extra rows copied from 2021-Q4 (which had 14 holdings) down to
# the 6 rows needed for 2022-Q1, shifting rows 8:15 up and out.
$newSheet.Range("A8:H15").Delete(-4162)

# Fund holdings for 2022-Q1.
$holdings = @(
    @("006608", "泓德研究优选混合", "21.75", "90.61", "2.71", "0.5894", 10),
    @("080012", "长盛电子信息产业混合", "6.40", "87.51", "3.70", "0.2368", 5),
    @("007305", "国联安新科技混合", "2.13", "81.55", "2.43", "0.0518", 10),
    @("003704", "光大保德信事件驱动灵活配置混合", "3.06", "23.55", "1.22", "0.0373", 7),
    @("009882", "华润元大核心动力混合A", "0.22", "68.63", "4.04", "0.0089", 7),
    @("009883", "华润元大核心动力混合C", "0.09", "68.63", "4.04", "0.0036", 7)
)

for ($i = 0; $i -lt $holdings.Length; $i++) {
    $r = $i + 2
    $rowData = $holdings[$i]

    # Fund code - force text so leading zeros survive, then strip the
    # quote-prefix styling so the cell matches plain inline-string cells.
    $newSheet.Cells.Item($r, 2).Value = "'" + $rowData[0]
    $newSheet.Cells.Item($r, 2).Style = "Normal"

    # Fund name - ordinary text, no special handling needed.
    $newSheet.Cells.Item($r, 3).Value = $rowData[1]

    # Numeric-looking text columns - force text so trailing zeros and exact
    # decimal formatting are preserved verbatim.
    $newSheet.Cells.Item($r, 4).Value = "'" + $rowData[2]
    $newSheet.Cells.Item($r, 4).Style = "Normal"
    $newSheet.Cells.Item($r, 5).Value = "'" + $rowData[3]
    $newSheet.Cells.Item($r, 5).Style = "Normal"
    $newSheet.Cells.Item($r, 6).Value = "'" + $rowData[4]
    $newSheet.Cells.Item($r, 6).Style = "Normal"
    $newSheet.Cells.Item($r, 7).Value = "'" + $rowData[5]
    $newSheet.Cells.Item($r, 7).Style = "Normal"

    # Rank column - true numeric value.
    $newSheet.Cells.Item($r, 8).Value = $rowData[6]
}

# Insert the 2022-Q1 summary row at the top of the "总计" sheet's data
# (row 2), pushing the existing quarters down. Re-fetch the sheet by name
# since the earlier sheet insert invalidates the cached reference.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Copy formatting (border/font/alignment) from the row below onto the new
# index cell so it matches the other "A" column cells, then set its value.
$totalSheet.Cells.Item(3, 1).Copy()
$totalSheet.Cells.Item(2, 1).PasteSpecial(-4122)
$totalSheet.Cells.Item(2, 1).Value = 0

# The date/count/value cells use plain (unstyled) cells like the other
# rows, so clear whatever formatting Insert() propagated into them.
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 6
$totalSheet.Cells.Item(2, 4).Value = 0.93

# The "A" column is a running 0-based row counter; renumber the rows that
# were pushed down (they kept their old 0,1,2 values from before the
# insert) to 1,2,3.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
